# Dashflux backup workbook enhancement:
#  - add 4 new lookup sheets (estabelecimentos, contas, categorias, cartoes)
#    extracted from the existing transaction sheets
#  - replace the old generic "Conta Padrão" account label with the real
#    bank name "ITAU" on the conta_corrente sheet

$wb = $excel.ActiveWorkbook

$contaCorrente = $wb.Worksheets.Item("conta_corrente")
$cartaoCredito = $wb.Worksheets.Item("cartao_credito")

# Canonical formatting sources already present on conta_corrente:
#   B1 -> bold header style (border + centered)
#   B2 -> plain/default body style
$headerFmt = $contaCorrente.Range("B1")
$bodyFmt = $contaCorrente.Range("B2")

# ---------------------------------------------------------------------
# 1) conta_corrente: the generic "Conta Padrão" account becomes "ITAU"
# ---------------------------------------------------------------------
$contaCorrente.Range("F2").Value = "ITAU"
$contaCorrente.Range("F3").Value = "ITAU"
$contaCorrente.Range("F4").Value = "ITAU"
$contaCorrente.Range("D3").Select()

$cartaoCredito.Range("F2").Select()

# ---------------------------------------------------------------------
# 2) estabelecimentos: distinct establishments used on conta_corrente
# ---------------------------------------------------------------------
$sEstabelecimentos = $wb.Worksheets.Add($null, $cartaoCredito)
$sEstabelecimentos.Name = "estabelecimentos"

$headerFmt.Copy()
$sEstabelecimentos.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$sEstabelecimentos.Range("A1").Value = "nome"
$sEstabelecimentos.Range("B1").Value = "descricao"

$bodyFmt.Copy()
$sEstabelecimentos.Range("A2:B3").PasteSpecial(-4122) | Out-Null
$sEstabelecimentos.Range("A2").Value = "MEU EMPREGO"
$sEstabelecimentos.Range("A3").Value = "PAGUE MENOS"

$sEstabelecimentos.Columns("A").ColumnWidth = 13
$sEstabelecimentos.Columns("B").ColumnWidth = 8.166666666666666

$sEstabelecimentos.Range("A2:A3").Select()

# ---------------------------------------------------------------------
# 3) contas: bank accounts referenced on conta_corrente
# ---------------------------------------------------------------------
$sContas = $wb.Worksheets.Add($null, $sEstabelecimentos)
$sContas.Name = "contas"

$headerFmt.Copy()
$sContas.Range("A1:F1").PasteSpecial(-4122) | Out-Null
$sContas.Range("A1").Value = "nome"
$sContas.Range("B1").Value = "descricao"
$sContas.Range("C1").Value = "se_banco"
$sContas.Range("D1").Value = "se_banco_nome"
$sContas.Range("E1").Value = "se_banco_agencia"
$sContas.Range("F1").Value = "se_banco_conta"

$bodyFmt.Copy()
$sContas.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$sContas.Range("A2").Value = "ITAU"
$sContas.Range("C2").Value = $true
$sContas.Range("D2").Value = "ITAU"

$sContas.Columns("C").ColumnWidth = 11.333333333333334
$sContas.Columns("D").ColumnWidth = 14.333333333333334
$sContas.Columns("E").ColumnWidth = 16
$sContas.Columns("F").ColumnWidth = 14

$sContas.Range("E2").Select()

# ---------------------------------------------------------------------
# 4) categorias: categories with an optional monthly goal ("meta")
# ---------------------------------------------------------------------
$sCategorias = $wb.Worksheets.Add($null, $sContas)
$sCategorias.Name = "categorias"

$headerFmt.Copy()
$sCategorias.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$sCategorias.Range("A1").Value = "nome"
$sCategorias.Range("B1").Value = "descricao"
$sCategorias.Range("C1").Value = "meta"

$bodyFmt.Copy()
$sCategorias.Range("A2:C2").PasteSpecial(-4122) | Out-Null
$sCategorias.Range("A2").Value = "SALÁRIO"
$sCategorias.Range("C2").Value = 1000

$sCategorias.Range("A3").Value = "FARMÁCIA"
$sCategorias.Range("C3").Value = 0

$sCategorias.Range("A4").Value = "SORVETE"

$sCategorias.Columns("A").ColumnWidth = 9.333333333333334
$sCategorias.Columns("B").ColumnWidth = 8.166666666666666

$sCategorias.Range("C4").Select()

# ---------------------------------------------------------------------
# 5) cartoes: credit cards referenced on cartao_credito
# ---------------------------------------------------------------------
$sCartoes = $wb.Worksheets.Add($null, $sCategorias)
$sCartoes.Name = "cartoes"

$headerFmt.Copy()
$sCartoes.Range("A1:E1").PasteSpecial(-4122) | Out-Null
$sCartoes.Range("A1").Value = "nome"
$sCartoes.Range("B1").Value = "descricao"
$sCartoes.Range("C1").Value = "bandeira"
$sCartoes.Range("D1").Value = "ultimos_4_digitos"
$sCartoes.Range("E1").Value = "dia_vencimento"

$bodyFmt.Copy()
$sCartoes.Range("B2:E2").PasteSpecial(-4122) | Out-Null
$sEstabelecimentos.Range("A2").Copy()
$sCartoes.Range("A2").PasteSpecial(-4122) | Out-Null
$sCartoes.Range("A2").Value = "CARTAO PRINCIPAL"
$sCartoes.Range("C2").Value = "OUTRO"
$sCartoes.Range("E2").Value = 16

$sCartoes.Columns("A").ColumnWidth = 17
$sCartoes.Columns("C").ColumnWidth = 7.666666666666667
$sCartoes.Columns("D").ColumnWidth = 15.5
$sCartoes.Columns("E").ColumnWidth = 14.166666666666666

$sCartoes.Range("E3").Select()

$contaCorrente.Select()
